$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are plain text in this sheet (e.g. "65.924.91",
# "0.998"), not real numbers. Force text format on the cells we are about
# to rewrite so numeric-looking strings like "0.998" or "6.20" aren't
# silently coerced to the Number type (and their formatting collapsed,
# e.g. "0.998" -> 0.998 displayed without trailing zero awareness, or
# "6.20" -> 6.2) by the COM layer's usual auto-detect-the-type behaviour.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.924.91"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "3.520.82"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "596.26"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("D6").Value = "143.48"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("D7").Value = "3.519.04"
$ws.Range("E7").Value = "  -1.36%  "
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("D11").Value = "7.68"
$ws.Range("E11").Value = "  -3.52%  "
$ws.Range("D12").Value = "0.405"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("D13").Value = "4.107.48"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "0.0000201"
$ws.Range("E14").Value = "  -3.54%  "
$ws.Range("D15").Value = "28.78"
$ws.Range("E15").Value = "  -4.50%  "
$ws.Range("D16").Value = "3.505.70"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "65.830.04"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("D19").Value = "10.96"
$ws.Range("E19").Value = "  -5.07%  "
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "14.40"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("D22").Value = "414.28"
$ws.Range("E22").Value = "  -4.68%  "
$ws.Range("D23").Value = "0.597"
$ws.Range("E23").Value = "  -2.43%  "
$ws.Range("D24").Value = "77.51"
$ws.Range("E24").Value = "  -2.99%  "
$ws.Range("D25").Value = "3.653.05"
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "0.0000116"
$ws.Range("E27").Value = "  -3.81%  "
$ws.Range("E28").Value = "  -2.02%  "
$ws.Range("D29").Value = "2.43"
$ws.Range("E29").Value = "  -3.11%  "
$ws.Range("D30").Value = "7.76"
$ws.Range("E30").Value = "  -3.74%  "
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "3.510.38"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("D33").Value = "0.155"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").Value = "24.34"
$ws.Range("E34").Value = "  -4.21%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "7.51"
$ws.Range("E36").Value = "  -4.82%  "
$ws.Range("E37").Value = "  -13.46%  "
$ws.Range("D38").Value = "175.45"
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("D39").Value = "5.23"
$ws.Range("E39").Value = "  -7.06%  "
$ws.Range("E40").Value = "  -8.38%  "
$ws.Range("D41").Value = "0.0821"
$ws.Range("E41").Value = "  -3.71%  "
$ws.Range("D42").Value = "5.05"
$ws.Range("E42").Value = "  -3.00%  "
$ws.Range("D43").Value = "0.860"
$ws.Range("E43").Value = "  -3.38%  "
$ws.Range("D44").Value = "45.31"
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("E45").Value = "  -8.36%  "
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  -4.76%  "
$ws.Range("D48").Value = "7.09"
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("D49").Value = "1.10"
$ws.Range("E49").Value = "  -8.04%  "
$ws.Range("D50").Value = "22.43"
$ws.Range("E50").Value = "  -5.01%  "
$ws.Range("D51").Value = "23.02"
$ws.Range("E51").Value = "  -8.48%  "
